$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-42 (C2:C42) -> 7295
$ws.Range("C2:C42").Value = 7295

# Rows 43-252 (C43:C252) -> 7293
$ws.Range("C43:C252").Value = 7293
